# WATER27.xlsx update
# The "DC" (dual-channel) SCAN functional columns are renamed/relabeled:
#   SCAN(DC)       -> SCAN        (column C)
#   SCAN           -> DC-SCAN     (column D)
#   SCAN(DC)-D3    -> SCAN-D3     (column E)
#   SCAN-D3        -> DC-SCAN-D3  (column F)
#   E_SCAN(DC)     -> E_SCAN      (column G)
#   E_SCAN         -> E_DC-SCAN   (column H)
#   E_SCAN(DC)-D3  -> E_SCAN-D3   (column I)
#   E_SCAN-D3      -> E_DC-SCAN-D3(column J)
# and the active cell selection moves from J13 to J2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "SCAN"
$ws.Range("D1").Value = "DC-SCAN"
$ws.Range("E1").Value = "SCAN-D3"
$ws.Range("F1").Value = "DC-SCAN-D3"
$ws.Range("G1").Value = "E_SCAN"
$ws.Range("H1").Value = "E_DC-SCAN"
$ws.Range("I1").Value = "E_SCAN-D3"
$ws.Range("J1").Value = "E_DC-SCAN-D3"

$ws.Range("J2").Select()
